$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep/store the value as literal text (matching the
    # source workbook's inline-string cells) instead of letting Excel infer
    # a numeric type for numeric-looking strings (e.g. "326.70", "-0.94%").
    $range.NumberFormat = "@"
    $range.Value = $value
    # Reset the visual style back to Normal/General now that the text type
    # has been locked in, so we do not leave a stray custom number format
    # on the cell.
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '326.70'
Set-TextValue $ws.Range('E2') '0.13%'
Set-TextValue $ws.Range('D3') '44.32'
Set-TextValue $ws.Range('E3') '-1.89%'
Set-TextValue $ws.Range('D4') '5.513'
Set-TextValue $ws.Range('E4') '-0.94%'
Set-TextValue $ws.Range('D5') '0.08020'
Set-TextValue $ws.Range('E5') '-0.89%'
Set-TextValue $ws.Range('D6') '1.989'
Set-TextValue $ws.Range('E6') '4.31%'
Set-TextValue $ws.Range('B7') 'BTSEToken'
Set-TextValue $ws.Range('C7') 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D7') '2.573'
Set-TextValue $ws.Range('E7') '-5.31%'
Set-TextValue $ws.Range('B8') 'MXToken'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D8') '0.9484'
Set-TextValue $ws.Range('E8') '-0.12%'
Set-TextValue $ws.Range('B9') 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D9') '0.1152'
Set-TextValue $ws.Range('E9') '-1.01%'
Set-TextValue $ws.Range('B10') 'WazirX'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D10') '0.1845'
Set-TextValue $ws.Range('E10') '-2.52%'
Set-TextValue $ws.Range('B11') 'MCDex'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D11') '12.18'
Set-TextValue $ws.Range('E11') '42.07%'
Set-TextValue $ws.Range('B12') 'MandalaExchangeToken'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D12') '0.09729'
Set-TextValue $ws.Range('E12') '-4.33%'
Set-TextValue $ws.Range('B13') 'BitrueCoin'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D13') '0.04607'
Set-TextValue $ws.Range('E13') '11.38%'
Set-TextValue $ws.Range('B14') 'BitMartToken'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D14') '0.1066'
Set-TextValue $ws.Range('E14') '0.21%'
Set-TextValue $ws.Range('B15') 'BitForexToken'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D15') '0.001266'
Set-TextValue $ws.Range('E15') '-0.60%'
Set-TextValue $ws.Range('B16') 'CoinExToken'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D16') '0.04081'
Set-TextValue $ws.Range('E16') '-4.60%'
Set-TextValue $ws.Range('B17') 'TigerCash'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D17') '0.005791'
Set-TextValue $ws.Range('E17') '-4.42%'
Set-TextValue $ws.Range('B18') 'LEO'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.366'
Set-TextValue $ws.Range('E18') '-6.98%'
Set-TextValue $ws.Range('B19') 'GateToken'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D19') '4.290'
Set-TextValue $ws.Range('E19') '-1.06%'
Set-TextValue $ws.Range('B20') 'BitpandaEcosystemToken'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D20') '0.3479'
Set-TextValue $ws.Range('E20') '-0.18%'
Set-TextValue $ws.Range('B21') 'ProBitToken'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D21') '0.1405'
Set-TextValue $ws.Range('E21') '2.30%'
Set-TextValue $ws.Range('B22') 'ZBToken'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D22') '0.2544'
Set-TextValue $ws.Range('E22') '-4.56%'
Set-TextValue $ws.Range('B23') 'BitKan'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D23') '0.001241'
Set-TextValue $ws.Range('E23') '0.25%'
Set-TextValue $ws.Range('B24') 'HotbitToken'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D24') '0.004313'
Set-TextValue $ws.Range('E24') '-6.11%'
Set-TextValue $ws.Range('E25') '-3.75%'
Set-TextValue $ws.Range('D26') '0.0003740'
Set-TextValue $ws.Range('E26') '-6.54%'
Set-TextValue $ws.Range('D38') '0.02569'
Set-TextValue $ws.Range('E38') '-4.47%'
Set-TextValue $ws.Range('D39') '0.05554'
Set-TextValue $ws.Range('E39') '-0.13%'
Set-TextValue $ws.Range('D40') '0.007530'
Set-TextValue $ws.Range('E40') '-2.27%'
Set-TextValue $ws.Range('D41') '0.1394'
Set-TextValue $ws.Range('E41') '-0.07%'
Set-TextValue $ws.Range('D42') '0.007607'
Set-TextValue $ws.Range('E42') '-32.89%'
Set-TextValue $ws.Range('E43') '-3.10%'
Set-TextValue $ws.Range('D44') '0.008507'
Set-TextValue $ws.Range('E44') '-2.10%'
Set-TextValue $ws.Range('D45') '0.00007108'
Set-TextValue $ws.Range('E45') '-0.32%'
Set-TextValue $ws.Range('D46') '0.00000000749'
Set-TextValue $ws.Range('E46') '-0.43%'
Set-TextValue $ws.Range('D47') '0.003526'
Set-TextValue $ws.Range('E47') '54.84%'
Set-TextValue $ws.Range('E48') '-4.33%'
Set-TextValue $ws.Range('D49') '0.00002098'
Set-TextValue $ws.Range('E49') '-0.43%'
Set-TextValue $ws.Range('D50') '0.0001998'
Set-TextValue $ws.Range('E50') '-0.43%'

Write-Output "Applied symbol list update"
